# Append the new game record (row 48) to the "Plan1" stats table and update
# the view/selection to match, mirroring the source workbook's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of match stats for columns A:U (in column order).
$newRow = @(47, 53, 47, 21, 12, 3, 4, 16, 7, 5, 5, 396, 369, 315, 291, 9, 3, 39, 19, 23, 8)

$targetRow = 48
for ($col = 1; $col -le $newRow.Length; $col++) {
    $ws.Cells.Item($targetRow, $col).Value = $newRow[$col - 1]
}

# Scroll the view back up to column M and leave the selection on the cell
# right after the newly appended row (matches the saved workbook's view).
[void]$ws.Activate()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 1
    $win.ScrollColumn = 13
}
[void]$ws.Range("V48").Select()
